$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values in columns A, B, E, F, G, H between row 2 and row 3
$cols = @("A", "B", "E", "F", "G", "H")

foreach ($col in $cols) {
    $r2 = "$col" + "2"
    $r3 = "$col" + "3"
    $v2 = $ws.Range($r2).Value2
    $v3 = $ws.Range($r3).Value2
    $ws.Range($r2).Value = $v3
    $ws.Range($r3).Value = $v2
}
